# Daily attendance processing - 2025-10-25 11:41:35
#
# The "Recorded By" column (G) lists the users who recorded/edited each
# attendance session, separated by ", ". For every row where that list is
# exactly two names and the second one is "System", swap the order so
# "System" is reported first (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # column G - "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[1] -eq "System") {
            $cell.Value = $parts[1] + ", " + $parts[0]
        }
    }
}
